# Update "Pais" (COVID stats) worksheet:
#  - swap the label order for three pairs of adjacent rows (data follows the new label,
#    per the source diff which reordered the shared-string table for these country pairs)
#  - refresh numeric statistics for a number of rows
#  - bump the "last updated" timestamp in A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 6 de Agosto de 2020 a las 09:45"

# --- Row 6: India (simple refresh, no reordering) -------------------------------------
$ws.Range("B6").Value = 1967700
$ws.Range("C6").Value = 4461
$ws.Range("D6").Value = 1329026
$ws.Range("E6").Value = 597902
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 33
$ws.Range("H6").Value = 40772

# --- Rows 37/38: Ucrania / Republica Dominicana swap places + refresh ----------------
$ws.Range("A37").Value = "Ucrania"
$ws.Range("B37").Value = 76808
$ws.Range("C37").Value = 1318
$ws.Range("D37").Value = 42524
$ws.Range("E37").Value = 32465
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 31
$ws.Range("H37").Value = 1819

$ws.Range("A38").Value = "Republica Dominicana"
$ws.Range("B38").Value = 75660
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 40122
$ws.Range("E38").Value = 34316
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 1222

# --- Row 54: Armenia (simple refresh) -------------------------------------------------
$ws.Range("B54").Value = 39819
$ws.Range("C54").Value = 233
$ws.Range("D54").Value = 31556
$ws.Range("E54").Value = 7491
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 772

# --- Rows 106/107: Hungria / Maldivas swap places + refresh --------------------------
$ws.Range("A106").Value = "Hungria"
$ws.Range("B106").Value = 4597
$ws.Range("C106").Value = 33
$ws.Range("D106").Value = 3463
$ws.Range("E106").Value = 534
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 600

$ws.Range("A107").Value = "Maldivas"
$ws.Range("B107").Value = 4594
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 2706
$ws.Range("E107").Value = 1869
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 19

# --- Row 110: Zimbabue (simple refresh) -----------------------------------------------
$ws.Range("B110").Value = 4339
$ws.Range("C110").Value = 118
$ws.Range("D110").Value = 1264
$ws.Range("E110").Value = 2991
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 3
$ws.Range("H110").Value = 84

# --- Row 141: Letonia (simple refresh; D141 unchanged) --------------------------------
$ws.Range("B141").Value = 1275
$ws.Range("C141").Value = 18
$ws.Range("E141").Value = 173

# --- Row 176: Camboya (simple refresh; B176/C176 unchanged) ---------------------------
$ws.Range("D176").Value = 210
$ws.Range("E176").Value = 33

# --- Rows 182/183: San Martin (Parte Holandesa) / Bermudas swap places + refresh -----
$ws.Range("A182").Value = "San Martin (Parte Holandesa)"
$ws.Range("B182").Value = 160
$ws.Range("C182").Value = 4
$ws.Range("D182").Value = 64
$ws.Range("E182").Value = 80
$ws.Range("F182").Value = 0
$ws.Range("G182").Value = 0
$ws.Range("H182").Value = 16

$ws.Range("A183").Value = "Bermudas"
$ws.Range("B183").Value = 157
$ws.Range("C183").Value = 0
$ws.Range("D183").Value = 144
$ws.Range("E183").Value = 4
$ws.Range("F183").Value = 0
$ws.Range("G183").Value = 0
$ws.Range("H183").Value = 9

# --- Rows 202/203: Timor Oriental / Santa Lucia swap places (stats identical) --------
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("A203").Value = "Santa Lucia"
